# Apply cryptos-list refresh (values scraped on 2023-03-28).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.010.17"
$ws.Range("E2").Value = "  -0.10%  "

# Row 3
$ws.Range("D3").Value = "1.743.44"
$ws.Range("E3").Value = "  +1.40%  "

# Row 4
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.50"
$ws.Range("E5").Value = "  -1.49%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.23%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4984"
$ws.Range("E7").Value = "  +7.92%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3587"
$ws.Range("E8").Value = "  +4.13%  "

# Row 9
$ws.Range("E9").Value = "  -0.33%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07273"
$ws.Range("E10").Value = "  -0.53%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.062"
$ws.Range("E11").Value = "  +0.97%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.27%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.30"
$ws.Range("E13").Value = "  +2.11%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.971"
$ws.Range("E14").Value = "  +1.64%  "

# Row 15
$ws.Range("D15").Value = "1.745.68"
$ws.Range("E15").Value = "  +1.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.871"
$ws.Range("E16").Value = "  -0.67%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.25"
$ws.Range("E17").Value = "  -2.45%  "

# Row 18
$ws.Range("E18").Value = "  -0.74%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06374"
$ws.Range("E19").Value = "  +0.64%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.31%  "

# Row 21
$ws.Range("E21").Value = "  +0.56%  "

# Row 22
$ws.Range("E22").Value = "  +1.55%  "

# Row 23
$ws.Range("D23").Value = "27.075.35"
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.34"
$ws.Range("E24").Value = "  +4.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.049"
$ws.Range("E25").Value = "  -4.85%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.32"
$ws.Range("E26").Value = "  -1.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.94"
$ws.Range("E27").Value = "  +2.48%  "

# Row 28
$ws.Range("D28").Value = "1.943.34"
$ws.Range("E28").Value = "  +1.43%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.189"
$ws.Range("E29").Value = "  +2.39%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.53"
$ws.Range("E30").Value = "  +1.02%  "

# Row 31
$ws.Range("E31").Value = "  +2.45%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09512"
$ws.Range("E32").Value = "  +4.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.577"
$ws.Range("E33").Value = "  -0.37%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.400"
$ws.Range("E34").Value = "  +1.27%  "

# Row 35
$ws.Range("E35").Value = "  +0.06%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05889"
$ws.Range("E36").Value = "  +0.63%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.09"
$ws.Range("E37").Value = "  -0.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.437"
$ws.Range("E38").Value = "  +2.32%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2004"
$ws.Range("E39").Value = "  +0.48%  "

# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.786"
$ws.Range("E40").Value = "  +0.26%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6038"
$ws.Range("E41").Value = "  +1.15%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.112"
$ws.Range("E42").Value = "  -1.30%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.591"
$ws.Range("E43").Value = "  +0.85%  "

# Row 44
$ws.Range("E44").Value = "  +2.22%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.592"
$ws.Range("E45").Value = "  -0.96%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5679"
$ws.Range("E46").Value = "  +0.73%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.40"
$ws.Range("E47").Value = "  +0.63%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.865"
$ws.Range("E48").Value = "  -0.08%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06679"
$ws.Range("E49").Value = "  +0.36%  "

# Row 50
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.101"
$ws.Range("E50").Value = "  +1.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  -0.22%  "
